$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data between row 2 and row 3 for columns A, I, K, P, Q, R and move AC value
$a2 = $ws.Range("A2").Value2
$a3 = $ws.Range("A3").Value2
$ws.Range("A2").Value2 = $a3
$ws.Range("A3").Value2 = $a2

$i2 = $ws.Range("I2").Value2
$i3 = $ws.Range("I3").Value2
$ws.Range("I2").Value2 = $i3
$ws.Range("I3").Value2 = $i2

$k2 = $ws.Range("K2").Value2
$k3 = $ws.Range("K3").Value2
$ws.Range("K2").Value2 = $k3
$ws.Range("K3").Value2 = $k2

$p2 = $ws.Range("P2").Value2
$p3 = $ws.Range("P3").Value2
$ws.Range("P2").Value2 = $p3
$ws.Range("P3").Value2 = $p2

$q2 = $ws.Range("Q2").Value2
$q3 = $ws.Range("Q3").Value2
$ws.Range("Q2").Value2 = $q3
$ws.Range("Q3").Value2 = $q2

$r2 = $ws.Range("R2").Value2
$r3 = $ws.Range("R3").Value2
$ws.Range("R2").Value2 = $r3
$ws.Range("R3").Value2 = $r2

$ac2 = $ws.Range("AC2").Value2
$ws.Range("AC3").Value2 = $ac2
$ws.Range("AC2").Value2 = ""
